# Add files via upload
# Fills in the "Normal Ogretim" (B/C) and "Ikinci Ogretim" (D/E) presenter
# names for the topic schedule on Sayfa1, and updates a couple of existing
# entries (Merve Nur Ates -> Merve Nur Ates(N.O.), Murat Arslan -> Murat
# Arslan, Abdullah Yildirim).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Row 5 - "3 Naive Bayes (Siniflandirma)"
$ws.Range("B5").Value = "Mustafa Görez, Mustafa Ercan (25 puan)"
$ws.Range("C5").Value = "Ramazan Dursun"
$ws.Range("E5").Value = "Ali Gökhan Çifter, Hamza Tayyip Şeker"

# Row 6 - "4 K-en yakin komsu (Siniflandirma)"
$ws.Range("C6").Value = "Beyza Demir"
$ws.Range("E6").Value = "Ömer Faruk Kıranşal, Enes Demir"

# Row 7 - "5 Yapay sinir aglari (Siniflandirma)"
$ws.Range("E7").Value = "Taner İğdirli, Yakup Çil"

# Row 8 - "6 Karar agaclari (Siniflandirma)"
$ws.Range("C8").Value = "Ali Yıldız"
$ws.Range("E8").Value = "Ceyhun Kayır, Selman Akaslan"

# Row 9 - "7 Destek vektor makinalari ve rastgele orman (Siniflandirma)"
$ws.Range("C9").Value = "Şule Nur Altun"
$ws.Range("E9").Value = "Soner Ateş, Yunus Emre Aras"

# Row 11 - "9 K-ortalama (Kumeleme)"
$ws.Range("C11").Value = "Furkan Emin Güven"
$ws.Range("E11").Value = "Umut Eroğlu, Muhammet Ali Kayran"

# Row 12 - "10 SLINK ve CLINK (Kumeleme)"
$ws.Range("C12").Value = "Arif Çelikkıran"
$ws.Range("E12").Value = "Özlem Donat, Ayşegül Çelik"

# Row 13 - "11 DBSCANS ve OPTICS (Kumeleme)"
$ws.Range("D13").Value = "Mehmet Onur Arslan, Büşra Polat"
$ws.Range("E13").Value = "Merve Nur Ateş(N.Ö.)"

# Row 14 - "12 Bagging, Boosting, Rastgele Altuzaylar (topluluk ogrenme)"
$ws.Range("E14").Value = "Sertaç Gülcan, Yahya Yıldız"

# Row 15 - "13 Siniflandirici performans olcumleri ..." (B15 text updated
# in place, same shared-string slot gains a second name)
$ws.Range("B15").Value = "Murat Arslan, Abdullah Yıldırım"

# Row 16 - "14 Regresyon"
$ws.Range("D16").Value = "Ömer Eşbah"

# Final selection left on E15, matching the last-edited cell
$ws.Range("E15").Select()
